# Update column F (dSF) values for specific rows per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = -6
$ws.Range("F10").Value = -1
$ws.Range("F15").Value = -4
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = -6
$ws.Range("F18").Value = -4
$ws.Range("F21").Value = -16
$ws.Range("F23").Value = -9
$ws.Range("F25").Value = -2
$ws.Range("F27").Value = 7
$ws.Range("F30").Value = -3
